$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
# Row 2 (order-date serial 5489)
$ws.Range("H2").Value = 1056.1
$ws.Range("J2").Value = 2149.6155
$ws.Range("L2").Value = 2149.6155
$ws.Range("N2").Value = -2375.6155
# Row 11 (order-date serial 5533)
$ws.Range("H11").Value = 47
$ws.Range("I11").Value = 47
$ws.Range("K11").Value = 47
$ws.Range("M11").Value = 93
# Row 32 (order-date serial 5484)
$ws.Range("H32").Value = 3747.5
$ws.Range("I32").Value = 1645
$ws.Range("J32").Value = 5850
$ws.Range("K32").Value = 1645
$ws.Range("L32").Value = 5850
$ws.Range("M32").Value = -1319
$ws.Range("N32").Value = -6502
# Row 41 (order-date serial 5478)
$ws.Range("H41").Value = 54
$ws.Range("I41").Value = 54
$ws.Range("K41").Value = 54
$ws.Range("M41").Value = 386
# Row 52 (order-date serial 4567)
$ws.Range("H52").Value = 2500
$ws.Range("I52").Value = 2500
$ws.Range("K52").Value = 7500
$ws.Range("M52").Value = -7340
# Row 53 (order-date serial 5479)
$ws.Range("H53").Value = 211.9375
$ws.Range("I53").Value = 121.22222
$ws.Range("K53").Value = 121.22222
$ws.Range("M53").Value = 515.77778
# Row 86 (order-date serial 12603)
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89 (order-date serial 12603)
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 99 (order-date serial 19883)
$ws.Range("H99").Value = 245.5
$ws.Range("I99").Value = 255.66667
$ws.Range("K99").Value = 767.00001
$ws.Range("M99").Value = 730.99999
# Row 101 (order-date serial 19884)
$ws.Range("H101").Value = 20000838
$ws.Range("I101").Value = 25000796
$ws.Range("K101").Value = 75002388
$ws.Range("M101").Value = -75000766

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
# Row 32 (order-date serial 44147)
$ws.Range("H32").Value = 2267338.5
$ws.Range("I32").Value = 2802842.8
$ws.Range("K32").Value = 2802842.8
$ws.Range("M32").Value = -2802555.8
# Row 38 (order-date serial 2260)
$ws.Range("H38").Value = 32500
$ws.Range("I38").Value = 40000
$ws.Range("J38").Value = 25000
$ws.Range("K38").Value = 40000
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = -39533
$ws.Range("N38").Value = -25934
# Row 97 (order-date serial 19941)
$ws.Range("H97").Value = 1110.1111
$ws.Range("I97").Value = 1011.75
$ws.Range("K97").Value = 1011.75
$ws.Range("M97").Value = -515.75
# Row 132 (order-date serial 43997)
$ws.Range("H132").Value = 2501.524
$ws.Range("I132").Value = 2501.524
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7504.572
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4974.572
$ws.Range("N132").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
# Row 7 (order-date serial 1602)
$ws.Range("H7").Value = 16931
$ws.Range("I7").Value = 529
$ws.Range("J7").Value = 33333
$ws.Range("K7").Value = 529
$ws.Range("L7").Value = 33333
$ws.Range("M7").Value = -416
$ws.Range("N7").Value = -33559
# Row 70 (order-date serial 15553)
$ws.Range("H70").Value = 300459
$ws.Range("J70").Value = 300459
$ws.Range("L70").Value = 300459
$ws.Range("N70").Value = -301045
# Row 73 (order-date serial 15553)
$ws.Range("H73").Value = 300459
$ws.Range("J73").Value = 300459
$ws.Range("L73").Value = 300459
$ws.Range("N73").Value = -302487
# Row 86 (order-date serial 12526)
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89 (order-date serial 12526)
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 134 (order-date serial 43998)
$ws.Range("H134").Value = 2473.75
$ws.Range("I134").Value = 2473.75
$ws.Range("K134").Value = 7421.25
$ws.Range("M134").Value = -4886.25

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
# Row 31 (order-date serial 44023)
$ws.Range("H31").Value = 1326.5
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 2000
$ws.Range("N31").Value = -2590
# Row 34 (order-date serial 44023)
$ws.Range("H34").Value = 1326.5
$ws.Range("J34").Value = 2000
$ws.Range("L34").Value = 2000
$ws.Range("N34").Value = -2404
# Row 35 (order-date serial 1627)
$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 1333.3334
$ws.Range("K35").Value = 1333.3334
$ws.Range("M35").Value = -1039.3334
# Row 62 (order-date serial 12580)
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65 (order-date serial 12580)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 94 (order-date serial 32934)
$ws.Range("H94").Value = 142548.88
$ws.Range("I94").Value = 372087
$ws.Range("J94").Value = 4826
$ws.Range("K94").Value = 372087
$ws.Range("L94").Value = 4826
$ws.Range("M94").Value = -371636
$ws.Range("N94").Value = -5728
# Row 99 (order-date serial 36198)
$ws.Range("H99").Value = 4201.125
$ws.Range("I99").Value = 3372.7144
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 3372.7144
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -1874.7144
$ws.Range("N99").Value = -12996
# Row 126 (order-date serial 36198)
$ws.Range("H126").Value = 4201.125
$ws.Range("I126").Value = 3372.7144
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 10118.1432
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -7648.143199999999
$ws.Range("N126").Value = -34940

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
# Row 8 (order-date serial 16734)
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("K8").Value = 600
$ws.Range("M8").Value = -461
# Row 11 (order-date serial 4745)
$ws.Range("H11").Value = 97.40000000000001
$ws.Range("I11").Value = 74
$ws.Range("J11").Value = 132.5
$ws.Range("K11").Value = 222
$ws.Range("L11").Value = 397.5
$ws.Range("M11").Value = -82
$ws.Range("N11").Value = -677.5
# Row 12 (order-date serial 4854)
$ws.Range("H12").Value = 305.86667
$ws.Range("I12").Value = 318.18182
$ws.Range("J12").Value = 272
$ws.Range("K12").Value = 954.54546
$ws.Range("L12").Value = 816
$ws.Range("M12").Value = -781.54546
$ws.Range("N12").Value = -1162
# Row 18 (order-date serial 36056)
$ws.Range("H18").Value = 783.2
$ws.Range("I18").Value = 783.2
$ws.Range("K18").Value = 2349.6
$ws.Range("M18").Value = -2180.6
# Row 38 (order-date serial 4860)
$ws.Range("H38").Value = 40
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 10
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 30
$ws.Range("M38").Value = 137
$ws.Range("N38").Value = -724
# Row 102 (order-date serial 19813)
$ws.Range("H102").Value = 4750
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 15000
$ws.Range("M102").Value = -12566
# Row 113 (order-date serial 27843)
$ws.Range("H113").Value = 1556.7273
$ws.Range("J113").Value = 1535.625
$ws.Range("L113").Value = 4606.875
$ws.Range("N113").Value = -8946.875
# Row 136 (order-date serial 44093)
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
# Row 138 (order-date serial 44105)
$ws.Range("H138").Value = 16668236
$ws.Range("I138").Value = 16668236
$ws.Range("K138").Value = 50004708
$ws.Range("M138").Value = -49999568

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
# Row 10 (order-date serial 4306)
$ws.Range("H10").Value = 12511000
$ws.Range("I10").Value = 25006000
$ws.Range("J10").Value = 15999.5
$ws.Range("K10").Value = 25006000
$ws.Range("L10").Value = 15999.5
$ws.Range("M10").Value = -25005831
$ws.Range("N10").Value = -16337.5
# Row 132 (order-date serial 44008)
$ws.Range("H132").Value = 2483.3333
$ws.Range("I132").Value = 2483.3333
$ws.Range("K132").Value = 7449.999899999999
$ws.Range("M132").Value = -4919.999899999999

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
# Row 7 (order-date serial 36249)
$ws.Range("H7").Value = 8659.200000000001
# Row 16 (order-date serial 5289)
$ws.Range("H16").Value = 829.8333
$ws.Range("I16").Value = 840.4
$ws.Range("K16").Value = 840.4
$ws.Range("M16").Value = -670.4
# Row 46 (order-date serial 5282)
$ws.Range("H46").Value = 3461
$ws.Range("I46").Value = 2047
$ws.Range("J46").Value = 4875
$ws.Range("K46").Value = 2047
$ws.Range("L46").Value = 4875
$ws.Range("M46").Value = -1859
$ws.Range("N46").Value = -5251
# Row 61 (order-date serial 27740)
$ws.Range("H61").Value = 2724.875
$ws.Range("I61").Value = 2685.5715
$ws.Range("K61").Value = 2685.5715
$ws.Range("M61").Value = -2483.5715
# Row 93 (order-date serial 19993)
$ws.Range("H93").Value = 3490.8
$ws.Range("I93").Value = 3051
$ws.Range("J93").Value = 5250
$ws.Range("K93").Value = 3051
$ws.Range("L93").Value = 5250
$ws.Range("M93").Value = -1803
$ws.Range("N93").Value = -7746
# Row 113 (order-date serial 27740)
$ws.Range("H113").Value = 2724.875
$ws.Range("I113").Value = 2685.5715
$ws.Range("K113").Value = 2685.5715
$ws.Range("M113").Value = -515.5715
# Row 126 (order-date serial 36249)
$ws.Range("H126").Value = 8659.200000000001
# Row 136 (order-date serial 44060)
$ws.Range("H136").Value = 2879.6
$ws.Range("I136").Value = 1959.6
$ws.Range("J136").Value = 3799.6
$ws.Range("K136").Value = 5878.799999999999
$ws.Range("L136").Value = 11398.8
$ws.Range("M136").Value = -3328.799999999999
$ws.Range("N136").Value = -16498.8

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
# Row 122 (order-date serial 36208)
$ws.Range("H122").Value = 2057.3333
$ws.Range("I122").Value = 2057.3333
$ws.Range("K122").Value = 6171.999899999999
$ws.Range("M122").Value = -3721.999899999999
